# db_schema clean-up: insert a new "type" column between "name" and
# "description" in the Tableau1 table on Sheet1, and fill it with each
# modality's data type (integer/string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects.Item("Tableau1")

# Shift "description" (and anything after it) one column to the right by
# inserting a worksheet column in its place. This leaves a blank column D
# that becomes the new "type" column, with "description" sliding into E.
$descCol = $table.ListColumns.Item("description")
$descCol.Range.EntireColumn.Insert()

# Grow the table range so it covers the newly inserted column.
$table.Resize($ws.Range("A1:E9"))

# Writing the header cells directly (rather than via ListColumn.Name)
# is what actually updates the ListObject's column metadata here.
$ws.Cells.Item(1, 4).Value = "type"
$ws.Cells.Item(1, 5).Value = "description"

# Per-row "type" values (row 1 is the header, data starts at row 2).
$types = @("integer", "string", "string", "integer", "string", "string", "integer", "integer")

for ($i = 0; $i -lt $types.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $types[$i]
}

# New "type" column is narrow; "description" keeps its original width
# (already carried over to column E by the column insert above).
$ws.Columns.Item(4).ColumnWidth = 6.33
